$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ParameterName (B2) to "testing" and ParameterID (A2) with a new generated id
$ws.Range("B2").Value = "testing"
$ws.Range("A2").Value = "36691e1f-7c63-40ef-ad40-5e0647a8492f"

# Update the active selection shown in the sheet view
$ws.Range("A8").Select()
